$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M31").Value = 739.1
$ws1.Range("L43").Value = 818.6900000000001
$ws1.Range("M43").Value = 2450.83
$ws1.Range("M44").Value = 3559.78
$ws1.Range("G51").Value = 68.03
$ws1.Range("G56").Value = "1 de 54"
$ws1.Range("M56").Value = "14 de 54"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F31").Value = 739.1
$ws2.Range("F43").Value = 3269.52
$ws2.Range("F44").Value = 5003.99
$ws2.Range("F51").Value = 1439.92
$ws2.Range("F56").Value = 82025.33

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
# Excel's ColumnWidth (chars) rounds to a pixel grid before being stored as the
# OOXML <col width> value; these inputs land exactly on width=13 and width=24.
$ws3.Columns.Item(4).ColumnWidth = 12.1
$ws3.Columns.Item(5).ColumnWidth = 23.1

$ws3.Range("D6").Value = 108.77
$ws3.Range("E6").Value = -1.950000000000003
$ws3.Range("F6").Value = 1.018255008425389

$ws3.Range("D15").Value = 16903.3
$ws3.Range("E15").Value = 3786.700000000001
$ws3.Range("F15").Value = 0.8169792170130498

$ws3.Range("D16").Value = 53150.7
$ws3.Range("E16").Value = 1570.530000000006
$ws3.Range("F16").Value = 0.9712994389928734

$ws3.Range("D19").Value = 88894.58
$ws3.Range("E19").Value = 16318.29000000001
$ws3.Range("F19").Value = 0.8449021493283094
